$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @{
    "A2"  = "0 - 4"
    "A3"  = "5 - 9"
    "A4"  = "10 - 14"
    "A5"  = "15 - 19"
    "A6"  = "20 - 24"
    "A7"  = "25 - 29"
    "A8"  = "30 - 34"
    "A9"  = "35 - 39"
    "A10" = "40 - 44"
    "A11" = "45 - 49"
    "A12" = "50 - 54"
    "A13" = "55 - 59"
    "A14" = "60 - 64"
    "A15" = "65 - 69"
    "A16" = "70 - 74"
    "A17" = "75 - 79"
    "A18" = "80 - 84"
    "A19" = "85 - 89"
    "A20" = "90 - 94"
    "A21" = "95 - 99"
    "A22" = "100+"
}

foreach ($addr in $values.Keys) {
    $ws.Range($addr).Value = $values[$addr]
}
